# The deck's theme (ppt/theme/theme1.xml, used by the slide master / all
# slides) is switched from the custom "Integral" palette to the stock
# Office 2016+ default palette ("Office Theme"): same dk1/lt1 (black/white),
# but dk2/lt2/accent1-6/hlink/folHlink change to the default Office colours.
#
# PowerPoint's ThemeColorScheme index order is:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# and .RGB is the usual COM "bgr-packed" integer (0x00BBGGRR).

function ToCOMColor($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

$officeThemeColors = @(
    (ToCOMColor 0x00 0x00 0x00),  # dk1      000000
    (ToCOMColor 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (ToCOMColor 0x44 0x54 0x6A),  # dk2      44546A
    (ToCOMColor 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (ToCOMColor 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (ToCOMColor 0xED 0x7D 0x31),  # accent2  ED7D31
    (ToCOMColor 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (ToCOMColor 0xFF 0xC0 0x00),  # accent4  FFC000
    (ToCOMColor 0x44 0x72 0xC4),  # accent5  4472C4
    (ToCOMColor 0x70 0xAD 0x47),  # accent6  70AD47
    (ToCOMColor 0x05 0x63 0xC1),  # hlink    0563C1
    (ToCOMColor 0x95 0x4F 0x72)   # folHlink 954F72
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
